$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: replace First/Last name and phone number with garbled inline-string values
$ws.Range("A4").Value = "JESSAy"
$ws.Range("B4").Value = "maerxhwelrl"
$ws.Range("D4").Value = "123-456-8678"

# Add a new row 12 duplicating the "JERSH / MERXWERLS / 123-111-9928" entry
# that already appears in rows 7-11
$ws.Range("A12").Value = "JERSH"
$ws.Range("B12").Value = "MERXWERLS"
$ws.Range("D12").Value = "123-111-9928"
